# Apply cryptos-list price/volume refresh per commit
# "Updated cryptos list on Fri Aug 25 10:46:26 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells whose new text is a plain decimal number (e.g. "1.004")
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a numeric value (losing the literal "1.004" text look).
$textPriceCells = @(
  "D4",
  "D5",
  "D6",
  "D9",
  "D11",
  "D12",
  "D15",
  "D17",
  "D20",
  "D21",
  "D23",
  "D25",
  "D26",
  "D27",
  "D28",
  "D30",
  "D31",
  "D32",
  "D33",
  "D35",
  "D36",
  "D37",
  "D38",
  "D39",
  "D40",
  "D41",
  "D42",
  "D44",
  "D47",
  "D49",
  "D51"
)
foreach ($addr in $textPriceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.201.21"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.658.89"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "217.14"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").Value = "0.5163"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").Value = "0.06270"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("D11").Value = "0.07755"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "4.483"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "1.653.42"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "1.886.11"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "0.5452"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "0.0₅8135"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "64.86"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "26.214.76"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "4.610"
$ws.Range("E20").Value = "  -3.22%  "
$ws.Range("D21").Value = "192.10"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "5.985"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "139.55"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "0.1220"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("D27").Value = "7.274"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "16.12"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "0.05926"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").Value = "1.273"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "3.544"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "3.274"
$ws.Range("E34").Value = "  -6.20%  "
$ws.Range("D35").Value = "0.9607"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").Value = "2.423"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "2.772"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").Value = "0.5671"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("D39").Value = "6.041"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "0.01590"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "0.8552"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "1.011.50"
$ws.Range("E43").Value = "  -7.36%  "
$ws.Range("D44").Value = "100.59"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "1.800.82"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "56.45"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "8.042"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "0.4208"
$ws.Range("E51").Value = "  -0.54%  "
